# Swap the "Office Theme" and "Integral" theme colour schemes between the
# two theme parts used by this deck.
#
# ppt/theme/theme1.xml is wired to the Notes Master and currently holds the
# stock "Office Theme" colours; ppt/theme/theme2.xml is wired to the Slide
# Master (and is therefore the theme actually seen on every slide) and
# currently holds the "Integral" colours. The edit swaps the two colour
# schemes in place so the Slide Master ends up with the Office Theme colours
# and the Notes Master ends up with the Integral colours (the file names /
# relationships themselves are untouched - only the RGB content of each
# theme's colour scheme changes).

$p = $ppt.ActivePresentation

# Slide master -> ppt/theme/theme2.xml (was "Integral", becomes "Office Theme")
$slideScheme = $p.SlideMaster.Theme.ThemeColorScheme

# Notes master -> ppt/theme/theme1.xml (was "Office Theme", becomes "Integral")
$notesScheme = $p.NotesMaster.Theme.ThemeColorScheme

# Colour scheme slots, in MsoThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

# Values that were originally in theme1.xml ("Office Theme")
$officeThemeRGB = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

# Values that were originally in theme2.xml ("Integral")
$integralRGB = @(0, 16777215, 5332805, 13754083, 3722137, 3646819, 2412774, 38860, 13611854, 10915127, 2465643, 158642)

for ($i = 1; $i -le 12; $i++) {
    # theme2.xml (Slide Master) now takes the Office Theme colours
    $slideScheme.Colors($i).RGB = $officeThemeRGB[$i - 1]
    # theme1.xml (Notes Master) now takes the Integral colours
    $notesScheme.Colors($i).RGB = $integralRGB[$i - 1]
}
